$d = $word.ActiveDocument

# The "Description" cell for test case 1 in the test-case table currently
# reads "Get all users". Expand the wording and mark the run with the
# (mangled) font-name / color formatting that appears in the target
# revision - this matches a stray "&quot" font name picked up from a
# copy/paste, rendered here as literal text: & + "quot".
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("Get all users", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng.Text = "Get all users messages in a JSON string"
    $rng.Font.Name = [char]38 + "quot"
    $rng.Font.Color = 0
}
